$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (the "Förändrad" date column) for rows 2 through 18:
# change the value from 45170 to 45174 for every row where it currently holds 45170.
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = 45174
    }
}
